$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Event")
$ws.Range("A2").Value = "Hello"
